$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits
#    right under the H1 title at the top of the document.
# -----------------------------------------------------------------
$metaParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $metaParaIndex = $i
        break
    }
}

if ($metaParaIndex -ne $null) {
    $metaPara = $d.Paragraphs.Item($metaParaIndex)
    $metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
    $metaRange.Delete()
}

# -----------------------------------------------------------------
# 2) Insert a new (bold) paragraph containing the page title right
#    before the final paragraph (the one that currently holds the
#    "Create a cartoon style image..." image-prompt text).
# -----------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Ra Deluxe Free | Ancient Egypt Themed Slot</w:t></w:r></w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($titleXml)

# -----------------------------------------------------------------
# 3) Swap the old image-prompt text for the meta-description copy
#    in the (now) very last paragraph, keeping its italic run
#    formatting intact.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon style image featuring an explorer inside a pyramid, surrounded by Egyptian symbols such as hieroglyphics and scarabs, with the Book of Ra symbol shining brightly in the background.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the treasures of ancient Egypt with Book of Ra Deluxe. Play for free and trigger exciting bonus features on this popular slot game.",
    2
)
